$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").Value = "Daniele Feltrinelli"
$ws.Range("B57").Value = "Stefano Tita | Clitoriders"
$ws.Range("C57").Value = "Federico  Manica | iMontagna"
$ws.Range("D57").Value = "Mattia Baldessarini | Shark Attack"
$ws.Range("E57").Value = "ANDREA ASTE | Pinguini Trentini"
$ws.Range("F57").Value = "Blendi Capa | FC. Stallions"
